$d = $word.ActiveDocument

# Split the "Xóa sp" bullet into its own paragraph, leaving a new empty
# paragraph (with the same ListParagraph/numbering) right after it.
$d.Content.Find.Execute("Xóa sp", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Xóa sp^p", 2)

# The new paragraph now holds the (moved) _GoBack bookmark; insert the new
# bullet text before it so the bookmark ends up after the run, as in the
# target revision.
$paras = $d.Paragraphs
$newPara = $paras.Item($paras.Count)
$newPara.Range.InsertBefore("TK2 thêm chức năng sửa sp")
